# Summary.xlsx RAD test-data regeneration
# Rewrites the execution log: a fresh Katalon test run (13/14 Oct 2024) replaces
# the previous run's timestamps, and the "Existing" sheet's Motor Fuel Tax row
# gains an MFLicNum column/value.
$wb = $excel.ActiveWorkbook

# --- Estimated: refresh run timestamps ---
$ws = $wb.Worksheets.Item("Estimated")
$vals = @("Sun Oct 13 11:10:24 EDT 2024", "Sun Oct 13 11:11:02 EDT 2024", "Sun Oct 13 11:11:36 EDT 2024", "Sun Oct 13 11:12:11 EDT 2024", "Sun Oct 13 11:12:46 EDT 2024", "Sun Oct 13 11:13:23 EDT 2024")
$row = 2
foreach ($v in $vals) {
    $ws.Cells.Item($row, 2).Value = $v
    $row = $row + 1
}

# --- Existing: add MFLicNum column (H) for the Motor Fuel Tax row ---
$ws = $wb.Worksheets.Item("Existing")
$ws.Columns.Item(8).ColumnWidth = 17.1666666667

$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "MFLicNum"

$ws.Range("C14").Copy($ws.Range("H14"))
$ws.Range("H14").Value = "Y"

# --- Existing: refresh run timestamps ---
$ws = $wb.Worksheets.Item("Existing")
$vals = @("Sun Oct 13 11:13:57 EDT 2024", "Sun Oct 13 11:14:32 EDT 2024", "Sun Oct 13 11:15:06 EDT 2024", "Sun Oct 13 11:15:41 EDT 2024", "Sun Oct 13 11:16:15 EDT 2024", "Sun Oct 13 11:16:50 EDT 2024", "Sun Oct 13 11:17:26 EDT 2024", "Sun Oct 13 11:18:03 EDT 2024", "Sun Oct 13 11:18:38 EDT 2024", "Sun Oct 13 11:19:14 EDT 2024", "Sun Oct 13 11:19:49 EDT 2024", "Sun Oct 13 11:20:23 EDT 2024", "Mon Oct 14 22:04:25 EDT 2024", "Sun Oct 13 11:22:11 EDT 2024", "Sun Oct 13 11:22:45 EDT 2024", "Sun Oct 13 11:23:21 EDT 2024", "Sun Oct 13 11:23:57 EDT 2024", "Sun Oct 13 11:24:31 EDT 2024")
$row = 2
foreach ($v in $vals) {
    $ws.Cells.Item($row, 2).Value = $v
    $row = $row + 1
}

# --- Extension: refresh run timestamps ---
$ws = $wb.Worksheets.Item("Extension")
$vals = @("Sun Oct 13 11:25:06 EDT 2024", "Sun Oct 13 11:25:40 EDT 2024", "Sun Oct 13 11:26:12 EDT 2024", "Sun Oct 13 11:26:45 EDT 2024", "Sun Oct 13 11:27:17 EDT 2024", "Sun Oct 13 11:27:49 EDT 2024")
$row = 2
foreach ($v in $vals) {
    $ws.Cells.Item($row, 2).Value = $v
    $row = $row + 1
}

# --- NewTaxReturn: refresh run timestamps ---
$ws = $wb.Worksheets.Item("NewTaxReturn")
$vals = @("Sun Oct 13 11:28:22 EDT 2024", "Sun Oct 13 11:28:56 EDT 2024", "Sun Oct 13 11:29:30 EDT 2024", "Sun Oct 13 11:30:05 EDT 2024", "Sun Oct 13 11:30:38 EDT 2024", "Sun Oct 13 11:31:12 EDT 2024", "Sun Oct 13 11:31:46 EDT 2024", "Sun Oct 13 11:32:21 EDT 2024", "Sun Oct 13 11:32:55 EDT 2024", "Sun Oct 13 11:33:29 EDT 2024", "Sun Oct 13 11:34:03 EDT 2024", "Sun Oct 13 11:34:37 EDT 2024", "Sun Oct 13 11:35:11 EDT 2024", "Sun Oct 13 11:35:45 EDT 2024", "Sun Oct 13 11:36:19 EDT 2024", "Sun Oct 13 11:36:53 EDT 2024", "Sun Oct 13 11:37:29 EDT 2024", "Sun Oct 13 11:38:05 EDT 2024", "Sun Oct 13 11:38:42 EDT 2024", "Sun Oct 13 11:39:16 EDT 2024", "Sun Oct 13 11:39:52 EDT 2024", "Sun Oct 13 11:40:28 EDT 2024", "Sun Oct 13 11:41:03 EDT 2024", "Sun Oct 13 11:41:39 EDT 2024", "Sun Oct 13 11:42:15 EDT 2024", "Sun Oct 13 11:42:50 EDT 2024", "Sun Oct 13 11:43:26 EDT 2024", "Sun Oct 13 11:44:02 EDT 2024", "Sun Oct 13 11:44:38 EDT 2024", "Sun Oct 13 11:45:13 EDT 2024", "Sun Oct 13 11:45:50 EDT 2024", "Sun Oct 13 11:46:24 EDT 2024", "Sun Oct 13 11:47:00 EDT 2024", "Sun Oct 13 11:47:37 EDT 2024", "Sun Oct 13 11:48:12 EDT 2024", "Sun Oct 13 11:48:48 EDT 2024", "Sun Oct 13 11:49:24 EDT 2024", "Sun Oct 13 11:49:59 EDT 2024", "Sun Oct 13 11:50:35 EDT 2024", "Sun Oct 13 11:51:11 EDT 2024", "Sun Oct 13 11:51:47 EDT 2024", "Sun Oct 13 11:52:22 EDT 2024", "Sun Oct 13 11:52:58 EDT 2024", "Sun Oct 13 11:53:32 EDT 2024", "Sun Oct 13 11:54:08 EDT 2024", "Sun Oct 13 11:54:43 EDT 2024", "Sun Oct 13 11:55:18 EDT 2024", "Sun Oct 13 11:55:54 EDT 2024", "Sun Oct 13 11:56:30 EDT 2024", "Sun Oct 13 11:57:05 EDT 2024", "Sun Oct 13 11:57:41 EDT 2024")
$row = 2
foreach ($v in $vals) {
    $ws.Cells.Item($row, 2).Value = $v
    $row = $row + 1
}

# --- Personal_IND: refresh run timestamps ---
$ws = $wb.Worksheets.Item("Personal_IND")
$vals = @("Sun Oct 13 11:58:17 EDT 2024", "Sun Oct 13 11:58:49 EDT 2024", "Sun Oct 13 11:59:22 EDT 2024", "Sun Oct 13 11:59:52 EDT 2024", "Sun Oct 13 12:00:22 EDT 2024", "Sun Oct 13 12:00:53 EDT 2024", "Sun Oct 13 12:01:23 EDT 2024", "Sun Oct 13 12:01:53 EDT 2024")
$row = 2
foreach ($v in $vals) {
    $ws.Cells.Item($row, 2).Value = $v
    $row = $row + 1
}

# --- Personal_JNT: refresh run timestamps ---
$ws = $wb.Worksheets.Item("Personal_JNT")
$vals = @("Sun Oct 13 12:02:24 EDT 2024", "Sun Oct 13 12:02:54 EDT 2024", "Sun Oct 13 12:03:25 EDT 2024", "Sun Oct 13 12:04:02 EDT 2024", "Sun Oct 13 12:04:38 EDT 2024")
$row = 2
foreach ($v in $vals) {
    $ws.Cells.Item($row, 2).Value = $v
    $row = $row + 1
}

# --- Personal_EL: refresh run timestamps ---
$ws = $wb.Worksheets.Item("Personal_EL")
$vals = @("Sun Oct 13 12:05:15 EDT 2024", "Sun Oct 13 12:05:53 EDT 2024")
$row = 2
foreach ($v in $vals) {
    $ws.Cells.Item($row, 2).Value = $v
    $row = $row + 1
}

# --- Selections (set non-active sheet first, active sheet last) ---
$wsNew = $wb.Worksheets.Item("NewTaxReturn")
$wsNew.Range("J1").Select()

$wsExisting = $wb.Worksheets.Item("Existing")
$wsExisting.Range("C2:C19").Select()
